$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.122.65"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.478.48"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.47"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.72"
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "4.069.32"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "3.477.17"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "64.094.24"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.17"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.98"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.40"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.49"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "3.616.94"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.156"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("D33").Value = "3.505.89"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.99"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.77"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.41"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.56"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.34"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.44"
$ws.Range("E46").Value = "  -6.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.73"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.899"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "2.332.44"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0254"
$ws.Range("E51").Value = "  -2.61%  "
